# Auto-generated edit script applying the cryptos.xlsx diff
# Updates Price (D) and Volume(1h) (E) columns, and swaps the
# VeChain/Bittensor rows (45/46) per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.889.46'
$ws.Range('E2').Value = '  -3.00%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.918.23'
$ws.Range('E3').Value = '  -3.70%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '584.98'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.06'
$ws.Range('E6').Value = '  -5.10%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -2.40%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.915.97'
$ws.Range('E9').Value = '  -3.72%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.93'
$ws.Range('E10').Value = '  +4.63%  '
$ws.Range('E11').Value = '  -4.37%  '
$ws.Range('E12').Value = '  -3.82%  '
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '33.58'
$ws.Range('E14').Value = '  -5.55%  '
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.399.23'
$ws.Range('E16').Value = '  -3.79%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.820.26'
$ws.Range('E17').Value = '  -3.06%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.76'
$ws.Range('E18').Value = '  -4.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '2.918.43'
$ws.Range('E19').Value = '  -3.74%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '430.06'
$ws.Range('E20').Value = '  -4.94%  '
$ws.Range('E21').Value = '  -4.53%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.681'
$ws.Range('E22').Value = '  -2.23%  '
$ws.Range('E23').Value = '  -4.81%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '80.33'
$ws.Range('E24').Value = '  -3.34%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '10.79'
$ws.Range('E25').Value = '  -2.76%  '
$ws.Range('E26').Value = '  -3.56%  '
$ws.Range('E27').Value = '  -3.02%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.21'
$ws.Range('E30').Value = '  -2.87%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.61'
$ws.Range('E31').Value = '  -3.03%  '
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.56'
$ws.Range('E33').Value = '  -3.65%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.106'
$ws.Range('E34').Value = '  -3.12%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0₃0871'
$ws.Range('E35').Value = '  +1.14%  '
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.64'
$ws.Range('E37').Value = '  -4.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.128'
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('E39').Value = '  -4.17%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '49.78'
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('E41').Value = '  -4.55%  '
$ws.Range('E42').Value = '  -4.90%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.293'
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '40.70'
$ws.Range('E44').Value = '  -2.91%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0349'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '376.28'
$ws.Range('E46').Value = '  -4.48%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.677.63'
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.48'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.73'
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('E51').Value = '  -1.66%  '
